$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 (radar_wiper / red): the RGB hex code changes from A3FF12 to FF0000
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = "FF0000"

# ---------------------------------------------------------------------------
# Row 7 (labels): colour name changes from "grey" to "red", hex from 808080
# to FF0000, and the D7 swatch fill changes from grey to red (matching D6).
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "red"
$ws.Range("E7").Value = "FF0000"
$ws.Range("D7").Interior.Color = 255

# ---------------------------------------------------------------------------
# Row 10 (small_text): colour name changes from "yellow" to "light grey",
# hex from FFFF00 to 808080, and the D10 swatch fill changes from yellow to
# grey.
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "light grey"
$ws.Range("E10").Value = "808080"
$ws.Range("D10").Interior.Color = 8421504

# ---------------------------------------------------------------------------
# Selection / view bookkeeping: the author ended up with the view scrolled
# back to A1 (no frozen/forced top-left cell) and the cursor on E14.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
